$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text formatting on D-column cells whose new values are numeric-looking,
# so Excel stores them as text (matching the original inlineStr string cells)
# instead of auto-converting to numbers.
$textCells = @("D5", "D7", "D8", "D9", "D10", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D22", "D23", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) { $ws.Range($addr).NumberFormat = "@" }

$ws.Range('D2').Value = '26.864.63'
$ws.Range('E2').Value = '  -0.92%  '

$ws.Range('D3').Value = '1.812.51'
$ws.Range('E3').Value = '  +0.42%  '

$ws.Range('E4').Value = '  -0.15%  '

$ws.Range('D5').Value = '309.38'
$ws.Range('E5').Value = '  +0.03%  '

$ws.Range('E6').Value = '  -0.11%  '

$ws.Range('D7').Value = '0.4316'
$ws.Range('E7').Value = '  +1.79%  '

$ws.Range('D8').Value = '0.3711'
$ws.Range('E8').Value = '  +2.59%  '

$ws.Range('D9').Value = '0.07263'
$ws.Range('E9').Value = '  +0.15%  '

$ws.Range('D10').Value = '0.8667'
$ws.Range('E10').Value = '  +2.44%  '

$ws.Range('B11').Value = 'WrappedEther'
$ws.Range('C11').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D11').Value = '2.023.94'
$ws.Range('E11').Value = '  +13.49%  '

$ws.Range('B12').Value = 'Solana'
$ws.Range('C12').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D12').Value = '20.96'
$ws.Range('E12').Value = '  +3.03%  '

$ws.Range('D13').Value = '6.635'
$ws.Range('E13').Value = '  +3.69%  '

$ws.Range('D14').Value = '5.359'
$ws.Range('E14').Value = '  +0.96%  '

$ws.Range('D15').Value = '0.06931'
$ws.Range('E15').Value = '  +1.87%  '

$ws.Range('D16').Value = '1.008'
$ws.Range('E16').Value = '  -0.01%  '

$ws.Range('D17').Value = '80.66'
$ws.Range('E17').Value = '  -0.83%  '

$ws.Range('D18').Value = '0.000008935'
$ws.Range('E18').Value = '  +1.81%  '

$ws.Range('D19').Value = '1.003'
$ws.Range('E19').Value = '  -0.04%  '

$ws.Range('D20').Value = '15.28'
$ws.Range('E20').Value = '  +1.31%  '

$ws.Range('D21').Value = '26.888.71'
$ws.Range('E21').Value = '  -0.59%  '

$ws.Range('D22').Value = '5.216'
$ws.Range('E22').Value = '  +2.37%  '

$ws.Range('D23').Value = '11.19'
$ws.Range('E23').Value = '  +0.78%  '

$ws.Range('D24').Value = '2.220.22'
$ws.Range('E24').Value = '  +7.71%  '

$ws.Range('D25').Value = '154.21'
$ws.Range('E25').Value = '  +0.63%  '

$ws.Range('D26').Value = '1.871'
$ws.Range('E26').Value = '  -4.24%  '

$ws.Range('D27').Value = '18.27'
$ws.Range('E27').Value = '  +0.31%  '

$ws.Range('D28').Value = '5.236'
$ws.Range('E28').Value = '  +4.03%  '

$ws.Range('D29').Value = '1.906'
$ws.Range('E29').Value = '  +14.72%  '

$ws.Range('D30').Value = '115.22'
$ws.Range('E30').Value = '  +1.07%  '

$ws.Range('D31').Value = '0.08938'
$ws.Range('E31').Value = '  -0.68%  '

$ws.Range('D32').Value = '0.7589'
$ws.Range('E32').Value = '  +2.99%  '

$ws.Range('D33').Value = '1.178'
$ws.Range('E33').Value = '  +7.27%  '

$ws.Range('D34').Value = '4.445'
$ws.Range('E34').Value = '  +1.37%  '

$ws.Range('D35').Value = '2.805'
$ws.Range('E35').Value = '  -2.59%  '

$ws.Range('D36').Value = '1.007'
$ws.Range('E36').Value = '  +0.32%  '

$ws.Range('D37').Value = '1.129'
$ws.Range('E37').Value = '  +4.49%  '

$ws.Range('E38').Value = '  +1.50%  '

$ws.Range('E39').Value = '  +0.73%  '

$ws.Range('D40').Value = '0.5091'
$ws.Range('E40').Value = '  +2.04%  '

$ws.Range('D41').Value = '0.1652'
$ws.Range('E41').Value = '  +1.04%  '

$ws.Range('D42').Value = '2.675'
$ws.Range('E42').Value = '  +2.37%  '

$ws.Range('D43').Value = '6.576'
$ws.Range('E43').Value = '  +10.14%  '

$ws.Range('D44').Value = '8.338'
$ws.Range('E44').Value = '  +2.62%  '

$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').Value = '10.47'
$ws.Range('E45').Value = '  +1.68%  '

$ws.Range('B46').Value = 'Quant'
$ws.Range('C46').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D46').Value = '106.99'
$ws.Range('E46').Value = '  +1.84%  '

$ws.Range('E47').Value = '  -0.10%  '

$ws.Range('B48').Value = 'Decentraland'
$ws.Range('C48').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D48').Value = '0.4586'
$ws.Range('E48').Value = '  +0.85%  '

$ws.Range('B49').Value = 'NEARProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D49').Value = '1.654'
$ws.Range('E49').Value = '  +2.90%  '

$ws.Range('D50').Value = '0.06297'
$ws.Range('E50').Value = '  -0.59%  '

$ws.Range('D51').Value = '1.818'
$ws.Range('E51').Value = '  +5.62%  '
